$d = $word.ActiveDocument

$d.Content.Find.Execute("996×9=8964", $true, $false, $false, $false, $false, $true, 1, $false, "471×6=2826", 2) | Out-Null
$d.Content.Find.Execute("150×7=1050", $true, $false, $false, $false, $false, $true, 1, $false, "570×9=5130", 2) | Out-Null
$d.Content.Find.Execute("436×8=3488", $true, $false, $false, $false, $false, $true, 1, $false, "143×5=715", 2) | Out-Null
$d.Content.Find.Execute("209×5=1045", $true, $false, $false, $false, $false, $true, 1, $false, "425×7=2975", 2) | Out-Null
$d.Content.Find.Execute("390×8=3120", $true, $false, $false, $false, $false, $true, 1, $false, "657×2=1314", 2) | Out-Null
$d.Content.Find.Execute("237×6=1422", $true, $false, $false, $false, $false, $true, 1, $false, "152×8=1216", 2) | Out-Null
$d.Content.Find.Execute("637×8=5096", $true, $false, $false, $false, $false, $true, 1, $false, "556×5=2780", 2) | Out-Null
$d.Content.Find.Execute("303×4=1212", $true, $false, $false, $false, $false, $true, 1, $false, "881×8=7048", 2) | Out-Null
$d.Content.Find.Execute("229×7=1603", $true, $false, $false, $false, $false, $true, 1, $false, "372×7=2604", 2) | Out-Null
$d.Content.Find.Execute("863×7=6041", $true, $false, $false, $false, $false, $true, 1, $false, "359×6=2154", 2) | Out-Null
$d.Content.Find.Execute("848×6=5088", $true, $false, $false, $false, $false, $true, 1, $false, "105×5=525", 2) | Out-Null
$d.Content.Find.Execute("306×2=612", $true, $false, $false, $false, $false, $true, 1, $false, "372×5=1860", 2) | Out-Null
$d.Content.Find.Execute("630×9=5670", $true, $false, $false, $false, $false, $true, 1, $false, "754×7=5278", 2) | Out-Null
$d.Content.Find.Execute("282×3=846", $true, $false, $false, $false, $false, $true, 1, $false, "686×9=6174", 2) | Out-Null
$d.Content.Find.Execute("305×7=2135", $true, $false, $false, $false, $false, $true, 1, $false, "779×2=1558", 2) | Out-Null
$d.Content.Find.Execute("966×8=7728", $true, $false, $false, $false, $false, $true, 1, $false, "563×6=3378", 2) | Out-Null
$d.Content.Find.Execute("561×5=2805", $true, $false, $false, $false, $false, $true, 1, $false, "221×5=1105", 2) | Out-Null
$d.Content.Find.Execute("550×5=2750", $true, $false, $false, $false, $false, $true, 1, $false, "788×2=1576", 2) | Out-Null
$d.Content.Find.Execute("779×8=6232", $true, $false, $false, $false, $false, $true, 1, $false, "787×8=6296", 2) | Out-Null
$d.Content.Find.Execute("159×4=636", $true, $false, $false, $false, $false, $true, 1, $false, "940×8=7520", 2) | Out-Null
$d.Content.Find.Execute("275×6=1650", $true, $false, $false, $false, $false, $true, 1, $false, "452×8=3616", 2) | Out-Null
$d.Content.Find.Execute("880×8=7040", $true, $false, $false, $false, $false, $true, 1, $false, "941×8=7528", 2) | Out-Null
$d.Content.Find.Execute("183×3=549", $true, $false, $false, $false, $false, $true, 1, $false, "908×8=7264", 2) | Out-Null
$d.Content.Find.Execute("879×4=3516", $true, $false, $false, $false, $false, $true, 1, $false, "410×4=1640", 2) | Out-Null
$d.Content.Find.Execute("737×4=2948", $true, $false, $false, $false, $false, $true, 1, $false, "175×4=700", 2) | Out-Null
